# Reflection-code update (01 Dec 2019): append a new data row (row 6) to the
# "TestData" sheet that duplicates the existing "DatadrivenTest" / AppURL
# row (row 5), including its hyperlink, and move the sheet selection down
# onto the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from row 5's formatting (style/number formats/fonts) by copying it
# down into the new row 6, then overwrite with the actual values.
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

$ws.Range("A6").Value = "DatadrivenTest"
$ws.Range("B6").Value = "Y"
$ws.Range("C6").Value = "Chrome"
$ws.Range("D6").Value = "https://twitter.com/login?lang=en"
$ws.Range("E6").Value = 9967887510
$ws.Range("F6").Value = "Tcs@1983"

# Give F6 the same mailto: hyperlink that F4/F5 already carry.
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:Tcs@1983")

# Hyperlinks.Add re-stamps the cell's style; restore the shared
# "Hyperlink" look (centered, underlined) that the rest of column F uses.
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F6").HorizontalAlignment = -4108

# Move the active selection onto the freshly added row, matching the
# workbook's new selection state.
$ws.Range("A6:F6").Select()
